{"js": "// Update the 25 division-fact cells in the practice table.\n// Each old expression is unique in the document, so we resolve every\n// search range first (against the ORIGINAL text) and only then perform\n// the text replacements. This matters because a couple of the new\n// values coincide with other cells' old values (e.g. \"88\u00f77=\" is both\n// an old value and a new value for a different cell) \u2014 replacing\n// sequentially without resolving ranges up front could clobber an\n// already-updated cell.\nconst pairs = [\n  [\"92\u00f78=\", \"71\u00f74=\"],\n  [\"55\u00f72=\", \"19\u00f78=\"],\n  [\"48\u00f75=\", \"36\u00f76=\"],\n  [\"70\u00f74=\", \"29\u00f76=\"],\n  [\"90\u00f74=\", \"80\u00f78=\"],\n  [\"58\u00f73=\", \"63\u00f79=\"],\n  [\"35\u00f78=\", \"68\u00f75=\"],\n  [\"20\u00f77=\", \"20\u00f78=\"],\n  [\"63\u00f73=\", \"37\u00f78=\"],\n  [\"38\u00f74=\", \"99\u00f79=\"],\n  [\"55\u00f73=\", \"69\u00f73=\"],\n  [\"12\u00f72=\", \"42\u00f74=\"],\n  [\"47\u00f72=\", \"30\u00f74=\"],\n  [\"50\u00f72=\", \"87\u00f78=\"],\n  [\"55\u00f78=\", \"89\u00f75=\"],\n  [\"94\u00f76=\", \"91\u00f77=\"],\n  [\"41\u00f75=\", \"88\u00f77=\"],\n  [\"41\u00f74=\", \"61\u00f72=\"],\n  [\"16\u00f77=\", \"34\u00f72=\"],\n  [\"35\u00f72=\", \"79\u00f75=\"],\n  [\"88\u00f77=\", \"86\u00f76=\"],\n  [\"72\u00f74=\", \"86\u00f73=\"],\n  [\"72\u00f72=\", \"59\u00f73=\"],\n  [\"79\u00f72=\", \"14\u00f76=\"],\n  [\"61\u00f77=\", \"92\u00f76=\"],\n];\n\nconst body = context.document.body;\n\n// 1) Resolve a search range for every old value before mutating anything.\nconst searchResults = pairs.map(([oldText]) =>\n  body.search(oldText, { matchCase: true, matchWholeWord: false })\n);\nawait context.sync();\n\n// 2) Now that all ranges are located against the pre-edit text, replace\n//    each one with its new value.\nfor (let i = 0; i < pairs.length; i++) {\n  const [, newText] = pairs[i];\n  const items = searchResults[i].items;\n  for (let j = 0; j < items.length; j++) {\n    items[j].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 25 division-fact cells in the practice table.\n#\n# We address every cell directly by (row, col) in the single table on the\n# page rather than doing a text-based Find/Replace. This sidesteps a\n# collision in the data: a couple of the NEW values are identical to\n# OTHER cells' OLD values (e.g. \"88\u00f77=\" is both an old cell value and the\n# new value of a different cell), so a naive sequential\n# Find.Execute(Replace:=wdReplaceAll) could re-match text that a prior\n# replacement had just written. Reading every old value first (keyed by\n# its cell position) and only then writing the new values avoids that\n# entirely.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$map = @{\n    \"92\u00f78=\" = \"71\u00f74=\"\n    \"55\u00f72=\" = \"19\u00f78=\"\n    \"48\u00f75=\" = \"36\u00f76=\"\n    \"70\u00f74=\" = \"29\u00f76=\"\n    \"90\u00f74=\" = \"80\u00f78=\"\n    \"58\u00f73=\" = \"63\u00f79=\"\n    \"35\u00f78=\" = \"68\u00f75=\"\n    \"20\u00f77=\" = \"20\u00f78=\"\n    \"63\u00f73=\" = \"37\u00f78=\"\n    \"38\u00f74=\" = \"99\u00f79=\"\n    \"55\u00f73=\" = \"69\u00f73=\"\n    \"12\u00f72=\" = \"42\u00f74=\"\n    \"47\u00f72=\" = \"30\u00f74=\"\n    \"50\u00f72=\" = \"87\u00f78=\"\n    \"55\u00f78=\" = \"89\u00f75=\"\n    \"94\u00f76=\" = \"91\u00f77=\"\n    \"41\u00f75=\" = \"88\u00f77=\"\n    \"41\u00f74=\" = \"61\u00f72=\"\n    \"16\u00f77=\" = \"34\u00f72=\"\n    \"35\u00f72=\" = \"79\u00f75=\"\n    \"88\u00f77=\" = \"86\u00f76=\"\n    \"72\u00f74=\" = \"86\u00f73=\"\n    \"72\u00f72=\" = \"59\u00f73=\"\n    \"79\u00f72=\" = \"14\u00f76=\"\n    \"61\u00f77=\" = \"92\u00f76=\"\n}\n\n$rows = $tbl.Rows.Count\n$cols = $tbl.Columns.Count\n\n# Pass 1: read every cell's current text and, for cells whose text is a\n# key in $map, remember which (row, col) needs which new value.\n$targets = New-Object System.Collections.ArrayList\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $txt = $cell.Range.Text\n        $txt = $txt.TrimEnd([char]13, [char]7)\n        if ($map.ContainsKey($txt)) {\n            [void]$targets.Add(@{ Row = $r; Col = $c; NewText = $map[$txt] })\n        }\n    }\n}\n\n# Pass 2: write the new values. Because the targets were resolved against\n# the original text in pass 1, this is safe even where new/old values\n# overlap.\nforeach ($t in $targets) {\n    $cell = $tbl.Cell($t.Row, $t.Col)\n    $cell.Range.Text = $t.NewText\n}\n\nWrite-Output \"replaced $($targets.Count) cells\"\n"}
